# Updated cryptos list on Mon Oct  2 11:59:31 UTC 2023 with GitHub Actions
#
# Applies the per-row Price (column D) / Volume(1h) (column E) refresh,
# plus the MXToken <-> RocketPoolETH row swap (rows 44 & 45).
#
# Note: several Price values are plain decimal-looking numbers (e.g. "7.90",
# "219.23"). Excel's COM Value setter auto-converts such strings to numeric
# values (dropping the trailing zero / exact text form), so those are written
# with a leading apostrophe to force them to stay text, matching the original
# inline-string cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.259.94"
$ws.Range("E2").Value = "  +3.76%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.729.10"
$ws.Range("E3").Value = "  +2.21%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'219.23"
$ws.Range("E5").Value = "  +1.22%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +0.08%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# Row 8 - Solana
$ws.Range("E8").Value = "  +3.91%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.68%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +1.34%  "

# Row 11 - TRON
$ws.Range("D11").Value = "'0.0896"
$ws.Range("E11").Value = "  +0.51%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.974.44"
$ws.Range("E12").Value = "  +2.32%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.736.11"
$ws.Range("E13").Value = "  +2.83%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  +0.83%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +1.67%  "

# Row 16 - Litecoin
$ws.Range("D16").Value = "'67.71"
$ws.Range("E16").Value = "  +0.23%  "

# Row 17 - WrappedBTC
$ws.Range("D17").Value = "28.266.01"
$ws.Range("E17").Value = "  +3.75%  "

# Row 18 - BitcoinCash
$ws.Range("D18").Value = "'245.41"

# Row 19 - ShibaInu
$ws.Range("D19").Value = "0.0₃0753"
$ws.Range("E19").Value = "  +0.92%  "

# Row 20 - Chainlink
$ws.Range("D20").Value = "'7.90"
$ws.Range("E20").Value = "  -3.08%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  -0.04%  "

# Row 22 - Uniswap
$ws.Range("D22").Value = "'4.64"
$ws.Range("E22").Value = "  +1.36%  "

# Row 23 - Avalanche
$ws.Range("E23").Value = "  +0.33%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.92%  "

# Row 25 - Monero
$ws.Range("D25").Value = "'149.60"
$ws.Range("E25").Value = "  +0.84%  "

# Row 26 - Cosmos
$ws.Range("D26").Value = "'7.50"
$ws.Range("E26").Value = "  +2.42%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "'16.64"
$ws.Range("E27").Value = "  +0.86%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  -0.07%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  -0.04%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  +2.63%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  +2.24%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +0.49%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("E33").Value = "  +0.36%  "

# Row 34 - Maker
$ws.Range("D34").Value = "1.487.28"
$ws.Range("E34").Value = "  -5.47%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -2.35%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  +2.30%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "'0.603"
$ws.Range("E37").Value = "  -0.47%  "

# Row 38 - HuobiToken
$ws.Range("E38").Value = "  +0.63%  "

# Row 39 - VeChain
$ws.Range("E39").Value = "  +0.95%  "

# Row 40 - WEMIXToken
$ws.Range("E40").Value = "  +0.33%  "

# Row 41 - Aave
$ws.Range("D41").Value = "'69.93"
$ws.Range("E41").Value = "  +0.52%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  -0.10%  "

# Row 43 - FraxShare
$ws.Range("E43").Value = "  -2.11%  "

# Row 44 & 45 - MXToken / RocketPoolETH swap places in the ranking
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.877.63"
$ws.Range("E44").Value = "  +1.96%  "

$ws.Range("B45").Value = "MXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D45").Value = "'2.29"
$ws.Range("E45").Value = "  +1.06%  "

# Row 46 - TrustWalletToken
$ws.Range("D46").Value = "'0.802"
$ws.Range("E46").Value = "  +1.67%  "

# Row 47 - RenderToken
$ws.Range("D47").Value = "'1.73"
$ws.Range("E47").Value = "  +7.19%  "

# Row 48 - BabyDogeCoin
$ws.Range("E48").Value = "  +4.23%  "

# Row 49 - Quant
$ws.Range("D49").Value = "'90.14"
$ws.Range("E49").Value = "  -1.15%  "

# Row 50 - EnergySwap
$ws.Range("D50").Value = "'8.20"
$ws.Range("E50").Value = "  -0.78%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -1.13%  "
